# "working on annotation fix"
# Renames Sheet1 -> Annotation, fixes up several annotation vocabulary
# strings in the shared-string table (entity->resource, #model->"# model",
# RDF->rdf, bare SBO/GO/FMA ids get their collection prefix prepended,
# Formula/Charge -> formula/charge), re-lays the column widths for the
# new 7-column layout, and moves the active selection to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rename the worksheet tab ---
$ws.Name = "Annotation"

# --- header: "entity" column renamed to "resource" ---
$ws.Range("F1").Value = "resource"

# --- "#model" -> "# model" ---
$ws.Range("A2").Value = "# model"

# --- qualifier column: RDF -> rdf (every data row) ---
$ws.Range("C3").Value = "rdf"
$ws.Range("C4").Value = "rdf"
$ws.Range("C7").Value = "rdf"
$ws.Range("C8").Value = "rdf"
$ws.Range("C9").Value = "rdf"
$ws.Range("C11").Value = "rdf"
$ws.Range("C12").Value = "rdf"
$ws.Range("C13").Value = "rdf"
$ws.Range("C15").Value = "rdf"
$ws.Range("C16").Value = "rdf"
$ws.Range("C17").Value = "rdf"
$ws.Range("C20").Value = "rdf"
$ws.Range("C21").Value = "rdf"
$ws.Range("C22").Value = "rdf"
$ws.Range("C25").Value = "rdf"
$ws.Range("C26").Value = "rdf"
$ws.Range("C27").Value = "rdf"
$ws.Range("C31").Value = "rdf"
$ws.Range("C32").Value = "rdf"

# --- entity column: prefix bare ontology ids with their collection ---
$ws.Range("F3").Value = "sbo/SBO:0000293"
$ws.Range("F4").Value = "go/GO:0008152"
$ws.Range("F7").Value = "sbo/SBO:0000290"
$ws.Range("F8").Value = "go/GO:0005615"
$ws.Range("F9").Value = "fma/FMA:70022"
$ws.Range("F11").Value = "sbo/SBO:0000290"
$ws.Range("F12").Value = "go/GO:0005886"
$ws.Range("F13").Value = "fma/FMA:63841"
$ws.Range("F15").Value = "sbo/SBO:0000290"
$ws.Range("F16").Value = "go/GO:0005623"
$ws.Range("F17").Value = "fma/FMA:68646"
$ws.Range("F20").Value = "sbo/SBO:0000027"
$ws.Range("F21").Value = "sbo/SBO:0000281"
$ws.Range("F22").Value = "sbo/SBO:0000186"
$ws.Range("F25").Value = "sbo/SBO:0000247"
$ws.Range("F26").Value = "sbo/SBO:0000247"
$ws.Range("F27").Value = "sbo/SBO:0000247"
$ws.Range("F31").Value = "sbo/SBO:0000185"
$ws.Range("F32").Value = "sbo/SBO:0000176"

# --- Formula / Charge rows lower-cased ---
$ws.Range("C28").Value = "formula"
$ws.Range("C29").Value = "charge"

# --- column widths for the new layout ---
# (values chosen so the saved width matches the target as closely as
# this engine's column-width rounding allows)
$ws.Columns.Item(1).ColumnWidth = 18.333333333333332
$ws.Columns.Item(2).ColumnWidth = 18.333333333333332
$ws.Columns.Item(3).ColumnWidth = 12.166666666666666
$ws.Columns.Item(4).ColumnWidth = 11
$ws.Columns.Item(5).ColumnWidth = 13
$ws.Columns.Item(6).ColumnWidth = 18.333333333333332
$ws.Columns.Item(7).ColumnWidth = 31.833333333333332

# --- move the active selection to A2 ---
$ws.Range("A2").Select()
